$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1. Refresh the "time_taken" query timestamps on the existing "data" sheet
#    (these literal values come from the re-run of the export that also
#    produced the new metadata tab below).
# ---------------------------------------------------------------------------
$dataSheet.Range("F2").Value = "2021-10-05 14:33:29.848899"
$dataSheet.Range("F3").Value = "2021-10-05 14:33:29.848909"
$dataSheet.Range("F4").Value = "2021-10-05 14:33:29.848914"
$dataSheet.Range("F5").Value = "2021-10-05 14:33:29.848917"
$dataSheet.Range("F6").Value = "2021-10-05 14:33:29.848921"
$dataSheet.Range("F7").Value = "2021-10-05 14:33:29.848951"
$dataSheet.Range("F8").Value = "2021-10-05 14:33:29.848961"
$dataSheet.Range("F9").Value = "2021-10-05 14:33:29.848966"
$dataSheet.Range("F10").Value = "2021-10-05 14:33:29.848970"
$dataSheet.Range("F11").Value = "2021-10-05 14:33:29.848974"
$dataSheet.Range("F12").Value = "2021-10-05 14:33:29.848977"
$dataSheet.Range("F13").Value = "2021-10-05 14:33:29.848980"
$dataSheet.Range("F14").Value = "2021-10-05 14:33:29.848982"

# ---------------------------------------------------------------------------
# 2. Add the new "metadata" tab, placed right after "data".
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $dataSheet)
$newSheet.Name = "metadata"

# Header row (row 1)
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Re-use the exact same bold/border/centered header style that "data" uses
# for its own header row, so no new style record is created.
$dataSheet.Range("B1").Copy()
$newSheet.Range("B1:G1").PasteSpecial(-4122)

# Row 2 - data values. A2 mirrors the "data" sheet's first index-column style.
$dataSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Cobblestone Malformations"
$newSheet.Range("C2").Value = 6

# D2 has to stay literal text "1.0" (not get auto-coerced into the number 1)
# while keeping the default/no style, exactly like every other unstyled data
# cell on this sheet. Stage the text in a scratch cell formatted as Text,
# reset D2's format to the sheet default, then copy just the value across.
$newSheet.Range("Z1").NumberFormat = "@"
$newSheet.Range("Z1").Value = "1.0"
$newSheet.Range("A1").Copy()
$newSheet.Range("D2").PasteSpecial(-4122)
$newSheet.Range("Z1").Copy()
$newSheet.Range("D2").PasteSpecial(-4163)
$newSheet.Range("Z1").Clear()

$newSheet.Range("E2").Value = "2021-01-18T20:40:48.077207Z"
$newSheet.Range("F2").Value = "2021-10-05 14:33:29.846088"
$newSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/6/?format=json"

$newSheet.Range("A1").Select()

# Keep "data" as the active/selected tab (it was the only - and therefore
# active - sheet before this edit).
$dataSheet.Activate()
$dataSheet.Range("A1").Select()
